$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.663.25'
$ws.Range('E2').Value = '  -1.43%  '

# Row 3
$ws.Range('D3').Value = '3.369.28'
$ws.Range('E3').Value = '  -0.56%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.22'
$ws.Range('E5').Value = '  -1.32%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.50'
$ws.Range('E6').Value = '  -0.92%  '

# Row 7
$ws.Range('E7').Value = '  +0.08%  '

# Row 8
$ws.Range('D8').Value = '3.367.53'
$ws.Range('E8').Value = '  -0.62%  '

# Row 10
$ws.Range('E10').Value = '  +1.39%  '

# Row 11
$ws.Range('E11').Value = '  -3.70%  '

# Row 12
$ws.Range('E12').Value = '  -2.34%  '

# Row 13
$ws.Range('D13').Value = '3.944.69'
$ws.Range('E13').Value = '  -0.55%  '

# Row 14
$ws.Range('E14').Value = '  +0.82%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.83'
$ws.Range('E15').Value = '  +0.47%  '

# Row 16
$ws.Range('D16').Value = '3.370.41'
$ws.Range('E16').Value = '  -0.49%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000170'
$ws.Range('E17').Value = '  -3.74%  '

# Row 18
$ws.Range('D18').Value = '60.808.18'

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.82'
$ws.Range('E19').Value = '  -0.32%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.65'
$ws.Range('E20').Value = '  -3.79%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.22'
$ws.Range('E21').Value = '  -2.68%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '372.19'

# Row 23
$ws.Range('D23').Value = '3.510.15'

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.545'
$ws.Range('E24').Value = '  -2.48%  '

# Row 25
$ws.Range('E25').Value = '  +0.04%  '

# Row 26
$ws.Range('E26').Value = '  -0.34%  '

# Row 27
$ws.Range('E27').Value = '  -1.75%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.175'
$ws.Range('E28').Value = '  +9.66%  '

# Row 29
$ws.Range('E29').Value = '  -3.54%  '

# Row 30
$ws.Range('E30').Value = '  +0.27%  '

# Row 32
$ws.Range('E32').Value = '  -1.84%  '

# Row 33
$ws.Range('E33').Value = '  -1.87%  '

# Row 34
$ws.Range('E34').Value = '  -0.01%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.27'
$ws.Range('E35').Value = '  -0.46%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.09'
$ws.Range('E36').Value = '  -4.45%  '

# Row 37
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.53'
$ws.Range('E37').Value = '  -2.04%  '

# Row 38
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.75'
$ws.Range('E38').Value = '  -1.11%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.24'
$ws.Range('E39').Value = '  -0.35%  '

# Row 40
$ws.Range('E40').Value = '  -3.48%  '

# Row 41
$ws.Range('E41').Value = '  -0.07%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.770'
$ws.Range('E42').Value = '  -1.23%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '25.05'
$ws.Range('E43').Value = '  -0.28%  '

# Row 44
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.64'
$ws.Range('E44').Value = '  +0.71%  '

# Row 45
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.67'
$ws.Range('E45').Value = '  -3.34%  '

# Row 46
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.29'
$ws.Range('E46').Value = '  -2.58%  '

# Row 47
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.17'
$ws.Range('E47').Value = '  -5.77%  '

# Row 48
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.515.79'
$ws.Range('E48').Value = '  +7.12%  '

# Row 49
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.21'
$ws.Range('E49').Value = '  +2.31%  '

# Row 50
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.74'
$ws.Range('E50').Value = '  -1.62%  '

# Row 51
$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.39'
$ws.Range('E51').Value = '  +3.45%  '
